$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows (26-37) extend the existing TimeSlice="UP" boundary table (rows 15-25)
# with four new Cset_CN groups, each repeated for years 2015/2020/2025, all
# boundary values set to 0 so that no production happens before 2030.
$rowsData = @(
  @{ Year = 2015; Attr = "CAP_BND"; Cset = "ERWINELCWIN3N" },
  @{ Year = 2020; Attr = "CAP_BND"; Cset = "ERWINELCWIN3N" },
  @{ Year = 2025; Attr = "CAP_BND"; Cset = "ERWINELCWIN3N" },
  @{ Year = 2015; Attr = "CAP_BND"; Cset = "ERWINELCWIN5N" },
  @{ Year = 2020; Attr = "CAP_BND"; Cset = "ERWINELCWIN5N" },
  @{ Year = 2025; Attr = "CAP_BND"; Cset = "ERWINELCWIN5N" },
  @{ Year = 2015; Attr = "ACT_BND"; Cset = "TB_ELCC_DKE_DKISLBH_01"; ExtraCol = "G" },
  @{ Year = 2020; Attr = "ACT_BND"; Cset = "TB_ELCC_DKE_DKISLBH_01"; ExtraCol = "G" },
  @{ Year = 2025; Attr = "ACT_BND"; Cset = "TB_ELCC_DKE_DKISLBH_01"; ExtraCol = "G" },
  @{ Year = 2015; Attr = "ACT_BND"; Cset = "TB_ELCC_DKW_DKISL1_01"; ExtraCol = "H" },
  @{ Year = 2020; Attr = "ACT_BND"; Cset = "TB_ELCC_DKW_DKISL1_01"; ExtraCol = "H" },
  @{ Year = 2025; Attr = "ACT_BND"; Cset = "TB_ELCC_DKW_DKISL1_01"; ExtraCol = "H" }
)

$destRow = 26
foreach ($rd in $rowsData) {
  $ws.Rows("25:25").Copy()
  $ws.Rows("$destRow`:$destRow").Insert(-4121)

  $ws.Range("D$destRow").Value = $rd.Year
  if ($rd.Attr -ne "CAP_BND") {
    $ws.Range("F$destRow").Value = $rd.Attr
  }
  if ($rd.ExtraCol) {
    $ws.Range("$($rd.ExtraCol)$destRow").Value = 0
  }
  $ws.Range("M$destRow").Value = $rd.Cset

  $destRow = $destRow + 1
}

$ws.Range("H38").Select()
